# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.247.27"
$ws.Cells.Item(2, 5).Value = "  +0.56%  "

$ws.Cells.Item(3, 4).Value = "1.857.29"
$ws.Cells.Item(3, 5).Value = "  +0.44%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "0.7063"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.79%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "238.00"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.08%  "

$ws.Cells.Item(7, 5).Value = "  +0.23%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.07977"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +2.58%  "

$ws.Cells.Item(9, 5).Value = "  -0.84%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "23.46"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.67%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.08173"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.77%  "

$ws.Cells.Item(12, 4).Value = "1.843.29"
$ws.Cells.Item(12, 5).Value = "  -0.35%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "5.188"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.45%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "0.7047"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -2.90%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "89.63"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.72%  "

$ws.Cells.Item(16, 4).Value = "29.267.59"
$ws.Cells.Item(16, 5).Value = "  +0.58%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.000007934"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.46%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "5.797"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.94%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "13.22"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.26%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "238.22"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.92%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "0.9989"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.04%  "

$ws.Cells.Item(22, 4).Value = "2.116.11"
$ws.Cells.Item(22, 5).Value = "  +0.57%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.26%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "7.474"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.67%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "162.83"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.14%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "8.864"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -1.27%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "0.1431"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.14%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "18.09"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +0.21%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.924"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -2.81%  "

$ws.Cells.Item(30, 5).Value = "  +2.06%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "1.475"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.68%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "4.367"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -2.70%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "4.018"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.29%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.05181"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.89%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.159"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -1.93%  "

$ws.Cells.Item(36, 5).Value = "  +1.31%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -1.92%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "2.651"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.32%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.01849"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.26%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "2.724"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.12%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.9356"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.35%  "

$ws.Cells.Item(42, 4).Value = "1.137.18"
$ws.Cells.Item(42, 5).Value = "  +3.75%  "

$ws.Cells.Item(43, 5).Value = "  -1.11%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.4258"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.16%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "70.30"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.47%  "

$ws.Cells.Item(46, 5).Value = "  +0.26%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "102.73"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.16%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.5307"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -4.10%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.759"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.54%  "

$ws.Cells.Item(50, 4).Value = "2.012.97"
$ws.Cells.Item(50, 5).Value = "  +0.77%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "9.164"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.06%  "

